$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the 2022 year header in column S, row 4, reusing the formatting
# of the adjacent 2021 column (R4).
$ws.Range("R4").Copy()
$ws.Range("S4").PasteSpecial(-4122)
$ws.Range("S4").Value = 2022

# Add the corresponding 2022 data value in column S, row 5, reusing the
# formatting of the adjacent 2021 column (R5).
$ws.Range("R5").Copy()
$ws.Range("S5").PasteSpecial(-4122)
$ws.Range("S5").Value = 76.1

# Update the selected/active cell shown in the saved worksheet view.
$ws.Range("P8").Select() | Out-Null
